$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ".demande" to each header cell in row 1 (columns B through J) that
# currently ends with ".jamais"
$cols = @("B","C","D","E","F","G","H","I","J")
foreach ($col in $cols) {
    $cell = $ws.Range("$col`1")
    $cell.Value = $cell.Value + ".demande"
}
